$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the "2020" (column L) year header and value columns with two more
# years, 2021 (M) and 2022 (N), copying the formatting of the existing 2020
# columns so the new cells look the same as their neighbours.
$ws.Range("L3").Copy($ws.Range("M3")) | Out-Null
$ws.Range("L3").Copy($ws.Range("N3")) | Out-Null
$ws.Range("M3").Value = 2021
$ws.Range("N3").Value = 2022

$ws.Range("L4").Copy($ws.Range("M4")) | Out-Null
$ws.Range("L4").Copy($ws.Range("N4")) | Out-Null
$ws.Range("M4").Value = 6.18
$ws.Range("N4").Value = 6.18

# Restore the active selection recorded in the saved view state.
$ws.Range("N15").Select() | Out-Null
